# Applies the changes described in the commit:
#  - Updated embedded documentation (comments)
#  - Added another row to the table
# (plus the renamed sheet / register metadata that goes with "testreg4")

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item(1)   # "registerinfo"
$wsReg  = $wb.Worksheets.Item(2)   # "testreg1" -> "testreg4"

# --- Rename the register worksheet -----------------------------------
$wsReg.Name = "testreg4"

# --- Update registerinfo values for the new register ------------------
$wsInfo.Range("B2").Value = "testreg4"
$wsInfo.Range("B3").Value = "4th test register"
$wsInfo.Range("B4").Value = "the description of the 4th test register"
$wsInfo.Range("B5").Value = "http://registry.it.csiro.au/sandbox/csiro/utils/testreg4"
$wsInfo.Range("B7").Value = "Simon Cox"

# Point the registry_location hyperlink at the new register's URL too
$wsInfo.Hyperlinks.Delete() | Out-Null
$wsInfo.Hyperlinks.Add($wsInfo.Range("B5"), "http://registry.it.csiro.au/sandbox/csiro/utils/testreg4") | Out-Null

# --- Update the column-header comments on the register sheet ----------
$wsReg.Range("B1").Comment.Text("A short label for the item, aka the 'term' which this concept definition is about") | Out-Null
$wsReg.Range("C1").Comment.Text("A description or definition of the item.") | Out-Null

# --- Add a new data row (item5) to the register table ------------------
$wsReg.Range("A6").Value = "item5"
$wsReg.Range("B6").Value = "Item 5"
$wsReg.Range("C6").Value = "Description of item 5"
$wsReg.Range("D6").Value = "i5"
$wsReg.Range("E6").Value = "to be noted regarding item 5: this is item 5 which is narrower than 4"
$wsReg.Range("F6").Value = "plucked from thin air"
$wsReg.Range("G6").Value = 4

# --- Selection state changes recorded in the saved file ----------------
$wsInfo.Range("A2").Select() | Out-Null
$wsReg.Activate() | Out-Null
$wsReg.Range("D10").Select() | Out-Null
